$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Files")

# Helper: write a value that must be stored as a shared-string (text) cell
# even though it looks numeric (e.g. "20161223", "10001"). A leading
# apostrophe forces Excel to keep it as text; ClearFormats() then drops the
# transient quote-prefix formatting so the cell ends up with no explicit
# style, matching how the sibling rows 2/3 store the same values.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.ClearFormats()
}

# Two more scanned images were found (image-00002, image-00003). Insert two
# new rows above the existing "dose_info.dcm" (unmapped) row - this pushes
# that row down to row 6 - and populate the new rows 4/5 the same way rows
# 2/3 are already populated.
$ws.Rows("4:5").Insert()

# Row 4: image-00002
$ws.Range("A4").Value = "Scans"
$ws.Range("B4").Value = "tests/fixtures/dose_info/billybob-10001/20161223/Scan/image-00002.dcm"
$ws.Range("C4").Value = "image-00002.dcm"
$ws.Range("D4").Value = "Y"
Set-TextValue $ws.Range("F4") "10001"
Set-TextValue $ws.Range("G4") "20161223"
$ws.Range("H4").Value = "Scan"
$ws.Range("I4").Value = "10001_CT1_6168"
$ws.Range("L4").Value = "CT"
$ws.Range("M4").Value = 6168
Set-TextValue $ws.Range("O4") "20161223"
$ws.Range("P4").Value = "CT1 abdomen"
$ws.Range("Q4").Value = "billybob"
Set-TextValue $ws.Range("R4") "10001"
$ws.Range("S4").Value = "Scan"
$ws.Range("T4").Value = "image-00002"

# Row 5: image-00003
$ws.Range("A5").Value = "Scans"
$ws.Range("B5").Value = "tests/fixtures/dose_info/billybob-10001/20161223/Scan/image-00003.dcm"
$ws.Range("C5").Value = "image-00003.dcm"
$ws.Range("D5").Value = "Y"
Set-TextValue $ws.Range("F5") "10001"
Set-TextValue $ws.Range("G5") "20161223"
$ws.Range("H5").Value = "Scan"
$ws.Range("I5").Value = "10001_CT1_6168"
$ws.Range("L5").Value = "CT"
$ws.Range("M5").Value = 6168
Set-TextValue $ws.Range("O5") "20161223"
$ws.Range("P5").Value = "CT1 abdomen"
$ws.Range("Q5").Value = "billybob"
Set-TextValue $ws.Range("R5") "10001"
$ws.Range("S5").Value = "Scan"
$ws.Range("T5").Value = "image-00003"

[void]$ws.Range("I5").Select()
